$d = $word.ActiveDocument

# --- 1. Remove the "On Pilgrimage..." heading and "By Dorothy Day" paragraphs,
#        including the bookmarkStart/bookmarkEnd that wrap them, and replace
#        with pandoc-title-block-style paragraphs (Title + Authors, multi-run).

$p1 = $d.Paragraphs.Item(1)
$p2 = $d.Paragraphs.Item(2)

# Clear paragraph 2's text but keep its paragraph mark.
$p2text = $d.Range($p2.Range.Start, $p2.Range.End - 1)
$p2text.Text = ""

# Delete paragraph 1 completely (text + paragraph mark). This causes the
# bookmarkEnd marker (which sat right after paragraph 1) to collapse down
# next to the bookmarkStart marker (which sits at document position 0).
$p1full = $d.Range($p1.Range.Start, $p1.Range.End)
$p1full.Delete()

# The two now-adjacent, now-empty bookmark markers sit at document position
# 0. A zero-length delete right at position 0 removes exactly one marker at
# a time, so do it twice to drop both bookmarkStart and bookmarkEnd.
$d.Range(0, 0).Delete()
$d.Range(0, 0).Delete()

# What remains is a single empty paragraph (the old "By Dorothy Day" slot)
# sitting right before "The Catholic Worker..." paragraph. Replace its
# content with the new Title paragraph + Authors paragraph, built from
# explicit OOXML so each word/space becomes its own run, matching the
# pandoc-generated markup.
$target = $d.Paragraphs.Item(1)
$newXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Title"/></w:pPr><w:r><w:t xml:space="preserve">On</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Pilgrimage</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">-</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">February</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">1957</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Authors"/></w:pPr><w:r><w:t xml:space="preserve">Dorothy</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Day</w:t></w:r></w:p>
'@
$target.Range.InsertXML($newXml)
